$d = $word.ActiveDocument

# Locate the run that currently reads "Рисунок  2 - рас" (the caption of
# figure 2, truncated mid-word "рас..."). We need to split the trailing
# "рас" into "Р" (plain/default run formatting) + "ас" (keeps the
# original Times New Roman / 28pt caption formatting), so the full
# caption becomes "Рисунок  2 - Расчеты для x = 6" once combined with the
# following run.
$needle = "Рисунок  2 - рас"
$found = $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text '$needle'"
}

$hit = $d.Content.Duplicate
$hit.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$tailLen = 3 # length of "рас"
$tailStart = $hit.Start + ($needle.Length - $tailLen)
$tailEnd = $hit.Start + $needle.Length

$tailRange = $d.Range($tailStart, $tailEnd)

# Replace the trailing "рас" with two runs: a plain "Р" run (default
# formatting, matching how Word represents a freshly-typed character
# with no explicit run properties besides rtl) followed by an "ас" run
# that keeps the exact original caption formatting (Times New Roman,
# 28pt/14pt, etc.).
$xmlFrag = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Р</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="none"/><w:shd w:fill="auto" w:val="clear"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">ас</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$tailRange.InsertXML($xmlFrag)
